$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all phone numbers in column C (rows 5 through 44) to the unified value
for ($r = 5; $r -le 44; $r++) {
    $ws.Cells.Item($r, 3).Value = "+5511961611974"
}

# Update the selection to match the target state
$ws.Range("C5:C44").Select()
